$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-CellText 2 4 '25.974.32'
Set-CellText 2 5 '  +0.51%  '
Set-CellText 3 4 '1.740.91'
Set-CellText 3 5 '  +0.30%  '
Set-CellText 4 4 '1.000'
Set-CellText 4 5 '  +0.03%  '
Set-CellText 5 4 '247.01'
Set-CellText 5 5 '  +4.38%  '
Set-CellText 6 4 '1.000'
Set-CellText 6 5 '  -0.05%  '
Set-CellText 7 4 '0.5027'
Set-CellText 7 5 '  -1.65%  '
Set-CellText 8 4 '0.2740'
Set-CellText 8 5 '  +0.64%  '
Set-CellText 9 4 '0.06188'
Set-CellText 9 5 '  +1.46%  '
Set-CellText 10 4 '1.745.90'
Set-CellText 10 5 '  +0.63%  '
Set-CellText 11 4 '0.07263'
Set-CellText 11 5 '  +1.36%  '
Set-CellText 12 4 '0.6539'
Set-CellText 12 5 '  +3.18%  '
Set-CellText 13 4 '15.17'
Set-CellText 13 5 '  +1.64%  '
Set-CellText 14 4 '4.695'
Set-CellText 14 5 '  +2.52%  '
Set-CellText 15 4 '77.55'
Set-CellText 15 5 '  +0.72%  '
Set-CellText 16 5 '  -0.05%  '
Set-CellText 17 4 '1.000'
Set-CellText 17 5 '  +0.00%  '
Set-CellText 18 4 '26.011.25'
Set-CellText 18 5 '  +0.68%  '
Set-CellText 19 4 '11.88'
Set-CellText 19 5 '  +1.54%  '
Set-CellText 20 4 '0.000006845'
Set-CellText 20 5 '  +2.30%  '
Set-CellText 21 4 '1.967.93'
Set-CellText 21 5 '  +0.19%  '
Set-CellText 22 4 '4.592'
Set-CellText 22 5 '  +8.35%  '
Set-CellText 23 4 '8.733'
Set-CellText 23 5 '  +1.53%  '
Set-CellText 24 4 '5.394'
Set-CellText 24 5 '  +3.42%  '
Set-CellText 25 4 '135.28'
Set-CellText 25 5 '  -2.63%  '
Set-CellText 26 4 '1.510'
Set-CellText 26 5 '  +0.07%  '
Set-CellText 27 4 '15.27'
Set-CellText 27 5 '  +1.05%  '
Set-CellText 28 4 '1.787'
Set-CellText 28 5 '  +2.24%  '
Set-CellText 29 4 '105.41'
Set-CellText 29 5 '  +0.22%  '
Set-CellText 30 4 '3.962'
Set-CellText 30 5 '  +1.37%  '
Set-CellText 31 4 '0.08151'
Set-CellText 31 5 '  -2.28%  '
Set-CellText 32 4 '3.707'
Set-CellText 32 5 '  +2.92%  '
Set-CellText 33 4 '0.04727'
Set-CellText 33 5 '  +3.96%  '
Set-CellText 34 4 '2.667'
Set-CellText 34 5 '  +0.67%  '
Set-CellText 35 4 '0.9979'
Set-CellText 35 5 '  +1.82%  '
Set-CellText 36 4 '0.6086'
Set-CellText 36 5 '  -1.72%  '
Set-CellText 37 4 '2.756'
Set-CellText 37 5 '  +2.71%  '
Set-CellText 38 4 '0.01621'
Set-CellText 38 5 '  +2.14%  '
Set-CellText 39 4 '1.933'
Set-CellText 39 5 '  +1.46%  '
Set-CellText 40 4 '0.9997'
Set-CellText 40 5 '  -0.05%  '
Set-CellText 41 4 '101.21'
Set-CellText 41 5 '  +3.95%  '
Set-CellText 42 4 '0.8073'
Set-CellText 42 5 '  +10.38%  '
Set-CellText 43 4 '0.3908'
Set-CellText 43 5 '  +2.11%  '
Set-CellText 44 4 '5.027'
Set-CellText 44 5 '  +2.05%  '
Set-CellText 45 4 '0.1170'
Set-CellText 45 5 '  +4.14%  '
Set-CellText 46 4 '6.362'
Set-CellText 46 5 '  +3.38%  '
Set-CellText 47 4 '55.83'
Set-CellText 47 5 '  +2.45%  '
Set-CellText 48 4 '0.05290'
Set-CellText 48 5 '  +0.44%  '
Set-CellText 49 4 '30.89'
Set-CellText 49 5 '  +1.82%  '
Set-CellText 50 4 '7.610'
Set-CellText 50 5 '  +0.94%  '
Set-CellText 51 4 '0.3469'
Set-CellText 51 5 '  +2.10%  '
